$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column E header
$ws.Range("E1").Value = "along"

# Data rows (country name, mean, CI_low, CI_high)
$data = @(
    @("<b>All</b>", 0.0107402699915947, -0.00289759134362644, 0.0243781313268159),
    @("<b>Europe</b>", 0.0165547590621099, -0.00318621188945287, 0.0362957300136726),
    @("France", -0.00494056137424157, -0.0512163291139684, 0.0413352063654853),
    @("Germany", 0.0481368196005593, 0.00722864651991834, 0.0890449926812002),
    @("Italy", 0.024627688715132, -0.022771040155234, 0.0720264175854981),
    @("Poland", -0.0129209851909682, -0.072844733673454, 0.0470027632915175),
    @("Spain", 0.0375599044291895, -0.0220015874743622, 0.0971213963327412),
    @("United Kingdom", 0.00518592478852574, -0.0428998635668623, 0.0532717131439138),
    @("Switzerland", -0.0135269516922421, -0.0757042568403408, 0.0486503534558566),
    @("Japan", 0.00955264618361334, -0.0173418224556559, 0.0364471148228826),
    @("Russia", -0.0182522013113218, -0.0630639428220745, 0.0265595401994309),
    @("Saudi Arabia", -0.00920322016048767, -0.0645346698358917, 0.0461282295149163),
    @("USA", 0.01779300719663, -0.00836073172973517, 0.0439467461229952)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = "info_solidarityTRUE"
    $row++
}
